$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (partial rich-text edits via Characters) ---
$ws.Range("A8").Characters(21, 2).Text = "51"
$ws.Range("C9").Characters(27, 10).Text = "12/19/2022"
$ws.Range("C9").Characters(48, 10).Text = "12/25/2022"

# --- Table cell updates ---
$ws.Range("D15").Value = 1
$ws.Range("C14").Copy($ws.Range("F15"))
$ws.Range("H15").Value = -100
$ws.Range("J15").Value = 13
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 85.714285714285
$ws.Range("C16").Value = 4
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = -36.842105263157
$ws.Range("I16").Value = 244
$ws.Range("J16").Value = 188
$ws.Range("K16").Value = 29.787234042553
$ws.Range("L16").Value = 70.629370629370
$ws.Range("M16").Value = 64.864864864864
$ws.Range("N16").Value = -71.327849588719
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 33.333333333333
$ws.Range("I17").Value = 177
$ws.Range("J17").Value = 172
$ws.Range("K17").Value = 2.906976744186
$ws.Range("L17").Value = 15.686274509803
$ws.Range("M17").Value = 88.297872340425
$ws.Range("N17").Value = -41.776315789473
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 28
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = 55.555555555555
$ws.Range("I18").Value = 416
$ws.Range("J18").Value = 238
$ws.Range("K18").Value = 74.789915966386
$ws.Range("L18").Value = 22.352941176470
$ws.Range("M18").Value = 137.714285714286
$ws.Range("N18").Value = -48.129675810473
$ws.Range("C19").Value = 24
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 140
$ws.Range("F19").Value = 107
$ws.Range("G19").Value = 83
$ws.Range("H19").Value = 28.915662650602
$ws.Range("I19").Value = 1364
$ws.Range("J19").Value = 862
$ws.Range("K19").Value = 58.236658932714
$ws.Range("L19").Value = 90.769230769230
$ws.Range("M19").Value = 30.526315789473
$ws.Range("N19").Value = -46.362563900904
$ws.Range("F14").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 1
$ws.Range("C14").Copy($ws.Range("D20"))
$ws.Range("E14").Copy($ws.Range("E20"))
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -40
$ws.Range("I20").Value = 58
$ws.Range("K20").Value = 5.454545454545
$ws.Range("L20").Value = 1.754385964912
$ws.Range("M20").Value = 28.888888888888
$ws.Range("N20").Value = -92.428198433420
$ws.Range("C21").Value = 36
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = 71.428571428571
$ws.Range("F21").Value = 159
$ws.Range("G21").Value = 134
$ws.Range("H21").Value = 18.656716417910
$ws.Range("I21").Value = 2275
$ws.Range("J21").Value = 1528
$ws.Range("K21").Value = 48.887434554973
$ws.Range("L21").Value = 60.663841807909
$ws.Range("M21").Value = 49.572649572649
$ws.Range("N21").Value = -56.904716802424
$ws.Range("F14").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("L14").Copy($ws.Range("E22"))
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 9
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 800
$ws.Range("I22").Value = 45
$ws.Range("J22").Value = 31
$ws.Range("K22").Value = 45.161290322580
$ws.Range("L22").Value = 55.172413793103
$ws.Range("M22").Value = -10
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 10
$ws.Range("G24").Value = 115
$ws.Range("H24").Value = 13.913043478260
$ws.Range("I24").Value = 1934
$ws.Range("J24").Value = 1297
$ws.Range("K24").Value = 49.113338473400
$ws.Range("L24").Value = 49.228395061728
$ws.Range("M24").Value = 33.471359558316
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -70
$ws.Range("G25").Value = 31
$ws.Range("H25").Value = -19.354838709677
$ws.Range("I25").Value = 408
$ws.Range("J25").Value = 320
$ws.Range("K25").Value = 27.5
$ws.Range("L25").Value = 65.182186234817
$ws.Range("M25").Value = 58.754863813229
$ws.Range("D26").Value = 1
$ws.Range("C14").Copy($ws.Range("F26"))
$ws.Range("H26").Value = -100
$ws.Range("J26").Value = 17
$ws.Range("K26").Value = 5.882352941176
$ws.Range("L26").Value = 100
$ws.Range("F14").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 2
$ws.Range("I27").Value = 78
$ws.Range("K27").Value = 14.705882352941
$ws.Range("L27").Value = 85.714285714285
